$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Roll the reporting periods forward by one fiscal year ---
# Row 8: period labels; Row 9: publish dates (text, not numeric)
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

$ws.Range("D9").Value = "1399-01-12 (7)"
$ws.Range("E9").Value = "1400-02-05 (7)"
$ws.Range("F9").Value = "1401-02-07 (9)"
$ws.Range("G9").Value = "1401-10-29 (7)"
$ws.Range("H9").Value = "1402-02-13 (2)"

# --- Shift all balance-sheet line items left one period, append new FY1401 figures ---
$ws.Range("D12").Value = 7605
$ws.Range("E12").Value = 68271
$ws.Range("F12").Value = 76456
$ws.Range("G12").Value = 118011
$ws.Range("H12").Value = 51199

$ws.Range("D13").Value = 265639
$ws.Range("E13").Value = 345656
$ws.Range("F13").Value = 396156
$ws.Range("G13").Value = 698811
$ws.Range("H13").Value = 616802

$ws.Range("D14").Value = 149658
$ws.Range("E14").Value = 202507
$ws.Range("F14").Value = 115674
$ws.Range("G14").Value = 95327
$ws.Range("H14").Value = 80455

$ws.Range("D15").Value = 361410
$ws.Range("E15").Value = 420406
$ws.Range("F15").Value = 429945
$ws.Range("G15").Value = 1151235
$ws.Range("H15").Value = 2014513

$ws.Range("D16").Value = 46667
$ws.Range("E16").Value = 43523
$ws.Range("F16").Value = 213548
$ws.Range("G16").Value = 200669
$ws.Range("H16").Value = 210733

$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 743
$ws.Range("F17").Value = 743
$ws.Range("G17").Value = 743
$ws.Range("H17").Value = 743

$ws.Range("D18").Value = 830979
$ws.Range("E18").Value = 1081106
$ws.Range("F18").Value = 1232522
$ws.Range("G18").Value = 2264796
$ws.Range("H18").Value = 2974445

$ws.Range("D19").Value = 14711
$ws.Range("E19").Value = 15128
$ws.Range("F19").Value = 16050
$ws.Range("G19").Value = 22992
$ws.Range("H19").Value = 36874

$ws.Range("D20").Value = 4938
$ws.Range("E20").Value = 4638
$ws.Range("F20").Value = 4638
$ws.Range("G20").Value = 4867
$ws.Range("H20").Value = 4867

$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0

$ws.Range("D22").Value = 361680
$ws.Range("E22").Value = 351037
$ws.Range("F22").Value = 354705
$ws.Range("G22").Value = 788445
$ws.Range("H22").Value = 1331920

$ws.Range("D23").Value = 9544
$ws.Range("E23").Value = 9468
$ws.Range("F23").Value = 9447
$ws.Range("G23").Value = 24439
$ws.Range("H23").Value = 26707

$ws.Range("D24").Value = '-'
$ws.Range("E24").Value = '-'
$ws.Range("F24").Value = '-'
$ws.Range("G24").Value = '-'
$ws.Range("H24").Value = '-'

$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

$ws.Range("D26").Value = 390873
$ws.Range("E26").Value = 380271
$ws.Range("F26").Value = 384840
$ws.Range("G26").Value = 840743
$ws.Range("H26").Value = 1400368

$ws.Range("D27").Value = 1221852
$ws.Range("E27").Value = 1461377
$ws.Range("F27").Value = 1617362
$ws.Range("G27").Value = 3105539
$ws.Range("H27").Value = 4374813

$ws.Range("D29").Value = 331850
$ws.Range("E29").Value = 390820
$ws.Range("F29").Value = 274738
$ws.Range("G29").Value = 744370
$ws.Range("H29").Value = 911813

$ws.Range("D30").Value = '-'
$ws.Range("E30").Value = '-'
$ws.Range("F30").Value = '-'
$ws.Range("G30").Value = '-'
$ws.Range("H30").Value = '-'

$ws.Range("D31").Value = 121793
$ws.Range("E31").Value = 75024
$ws.Range("F31").Value = 165289
$ws.Range("G31").Value = 231734
$ws.Range("H31").Value = 309502

$ws.Range("D32").Value = 154368
$ws.Range("E32").Value = 190303
$ws.Range("F32").Value = 117312
$ws.Range("G32").Value = 208189
$ws.Range("H32").Value = 267029

$ws.Range("D33").Value = 61652
$ws.Range("E33").Value = 7074
$ws.Range("F33").Value = 30624
$ws.Range("G33").Value = 30468
$ws.Range("H33").Value = 19396

$ws.Range("D34").Value = 63614
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 143616
$ws.Range("H34").Value = 128528

$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0

$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0

$ws.Range("D37").Value = 733277
$ws.Range("E37").Value = 663221
$ws.Range("F37").Value = 587963
$ws.Range("G37").Value = 1358377
$ws.Range("H37").Value = 1636268

$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0

$ws.Range("D39").Value = '-'
$ws.Range("E39").Value = '-'
$ws.Range("F39").Value = '-'
$ws.Range("G39").Value = '-'
$ws.Range("H39").Value = '-'

$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0

$ws.Range("D41").Value = 39600
$ws.Range("E41").Value = 38718
$ws.Range("F41").Value = 48818
$ws.Range("G41").Value = 69770
$ws.Range("H41").Value = 131294

$ws.Range("D42").Value = 39600
$ws.Range("E42").Value = 38718
$ws.Range("F42").Value = 48818
$ws.Range("G42").Value = 69770
$ws.Range("H42").Value = 131294

$ws.Range("D43").Value = 772877
$ws.Range("E43").Value = 701939
$ws.Range("F43").Value = 636781
$ws.Range("G43").Value = 1428147
$ws.Range("H43").Value = 1767562

$ws.Range("D45").Value = 125000
$ws.Range("E45").Value = 125000
$ws.Range("F45").Value = 125000
$ws.Range("G45").Value = 125000
$ws.Range("H45").Value = 200000

$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0

$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0

$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = -128891
$ws.Range("G48").Value = -74120
$ws.Range("H48").Value = -119179

$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0

$ws.Range("D50").Value = 12500
$ws.Range("E50").Value = 12500
$ws.Range("F50").Value = 12500
$ws.Range("G50").Value = 12500
$ws.Range("H50").Value = 20000

$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0

$ws.Range("D52").Value = '-'
$ws.Range("E52").Value = '-'
$ws.Range("F52").Value = '-'
$ws.Range("G52").Value = '-'
$ws.Range("H52").Value = '-'

$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0

$ws.Range("D54").Value = '-'
$ws.Range("E54").Value = '-'
$ws.Range("F54").Value = '-'
$ws.Range("G54").Value = '-'
$ws.Range("H54").Value = '-'

$ws.Range("D55").Value = 0
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 0

$ws.Range("D56").Value = 311475
$ws.Range("E56").Value = 621938
$ws.Range("F56").Value = 971972
$ws.Range("G56").Value = 1614012
$ws.Range("H56").Value = 2506430

$ws.Range("D57").Value = 448975
$ws.Range("E57").Value = 759438
$ws.Range("F57").Value = 980581
$ws.Range("G57").Value = 1677392
$ws.Range("H57").Value = 2607251

$ws.Range("D58").Value = 1221852
$ws.Range("E58").Value = 1461377
$ws.Range("F58").Value = 1617362
$ws.Range("G58").Value = 3105539
$ws.Range("H58").Value = 4374813
